# docs: rework some content
#
# - Journal de travail: fill in rows 27-29 (new work-log entries for 5/1/2023)
#   and let the totals formula (C30) recalc.
# - Compte rendu: reword the last entry (row 9) and resize things a bit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Journal de travail")
$ws2 = $wb.Worksheets.Item("Compte rendu")

# ---------------------------------------------------------------------------
# Sheet "Journal de travail"
# ---------------------------------------------------------------------------

# Row 27 - Implémentation: DB integration for module persistence + JWT auth
$ws1.Range("A27").Value = 45047
$ws1.Range("B27").Value = "Implémentation"
$ws1.Range("C27").Value = 6
$ws1.Range("D27").Value = "Intégration de la DB pour la persistance des modules et authentification des utilisateurs (JWT)"

# Row 28 - Rédaction: report writing
$ws1.Range("A28").Value = 45047
$ws1.Range("B28").Value = "Rédaction"
$ws1.Range("C28").Value = 2
$ws1.Range("D28").Value = "Rapport"

# Row 29 - Réunion avec professeur (new row, needs the same date style as
# the rows above it - copy it over rather than setting NumberFormat so we
# reuse the existing style record instead of minting a new one)
$ws1.Range("A28").Copy()
$ws1.Range("A29").PasteSpecial(-4122)
$ws1.Range("A29").Value = 45047
$ws1.Range("B29").Value = "Réunion"
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = "Réunion avec professeur"

# ---------------------------------------------------------------------------
# Sheet "Compte rendu"
# ---------------------------------------------------------------------------

$ws2.Range("C9").Value = "Retour sur la discussion avec l'entreprise et sur le travail effectué, conseils sur la rédaction du rapport"

# Widen the "Points importants" column and let row 8 shrink back down now
# that more of its text fits on fewer lines.
$ws2.Columns.Item(3).ColumnWidth = 85.66666666666667
$ws2.Rows.Item(8).RowHeight = 34

# ---------------------------------------------------------------------------
# Selections - update the cursor position on each sheet, making sure the
# journal sheet ends up re-activated last so it keeps being the visible tab.
# ---------------------------------------------------------------------------

$ws2.Activate() | Out-Null
$ws2.Range("C12").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D28").Select() | Out-Null
